$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right marks 5 -> 4, Wrong marks -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 110 -> 88, Wrong total -1 -> -2, Max text "110 / 140" -> "86 / 112"
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "86 / 112"
